$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1 & 2: paragraph "1.3. Esfuerzo y dedicación al proyecto(38%):(EXPLICAR)"
# Word's grammar checker had split this into extra runs around proofErr
# (gramStart/gramEnd) marks. Re-typing/normalizing the text (via a
# scoped Find & Replace with identical text) causes Word to re-tokenize the
# affected runs, merging "...proyecto(38" + "%)" into one run and "(" +
# "EXPLICAR)" into another, and drops the now-stale proofErr marks.
# The Find/Replace is scoped to just this paragraph's Range so the identical
# "(EXPLICAR)" text appearing in other paragraphs is left untouched.
# ---------------------------------------------------------------------------
$pEsfuerzo = $d.Paragraphs(5)

$r1 = $d.Range($pEsfuerzo.Range.Start, $pEsfuerzo.Range.End)
$r1.Find.Execute("1.3. Esfuerzo y dedicación al proyecto(38%)", $true, $false, $false, $false, $false, $true, 1, $false, "1.3. Esfuerzo y dedicación al proyecto(38%)", 2) | Out-Null

$r2 = $d.Range($pEsfuerzo.Range.Start, $pEsfuerzo.Range.End)
$r2.Find.Execute("(EXPLICAR)", $true, $false, $false, $false, $false, $true, 1, $false, "(EXPLICAR)", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 3: paragraph "2ª persona():" gets the name "Alberto" typed into the
# empty parentheses, matching the pattern already used for the other team
# members ("1ª persona (Augusto):", "4ª persona(Adrian):", ...).
# ---------------------------------------------------------------------------
$pPersona2 = $d.Paragraphs(18)

# Re-type the paragraph text (minus the end-of-paragraph mark) so the
# grammar-check proofErr marks bracketing "persona(" get cleared and
# "2ª " + "persona(" collapse into a single run, just like typing the name
# would invalidate the stale grammar check span.
$rFull = $d.Range($pPersona2.Range.Start, $pPersona2.Range.End - 1)
$rFull.Text = "2ª persona(Alberto):"

# Locate the newly-typed name and force it into its own run (mirrors how a
# real edit leaves the typed word as a distinct run from its neighbours) by
# toggling a character property on and back off.
$pPersona2 = $d.Paragraphs(18)
$rName = $d.Range($pPersona2.Range.Start, $pPersona2.Range.End)
$rName.Find.Execute("Alberto", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null
$rName.Bold = 1
$rName.Bold = 0

# The "_GoBack" bookmark (Word's "last edit location" marker) moves to sit
# right after the newly-inserted name - adding it here automatically removes
# it from its previous location (a bookmark name is unique per document).
$rGoBack = $d.Range($rName.End, $rName.End)
$d.Bookmarks.Add("_GoBack", $rGoBack) | Out-Null

# ---------------------------------------------------------------------------
# The paragraph "Puntos: 2.1" + "." used to be split by the old _GoBack
# bookmark; now that the bookmark has moved away, re-normalize it back into
# a single run "Puntos: 2.1."
# ---------------------------------------------------------------------------
$pPuntos21 = $d.Paragraphs(21)
$r3 = $d.Range($pPuntos21.Range.Start, $pPuntos21.Range.End)
$r3.Find.Execute("Puntos: 2.1.", $true, $false, $false, $false, $false, $true, 1, $false, "Puntos: 2.1.", 2) | Out-Null
